# Append two new tracker rows (64 and 65) to Sheet1, matching the pattern
# of the existing rows: GoalID, GoalName, Date, Progress, Percentage, Change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Cells.Item(63, 3).NumberFormat()

# Row 64: G1 / Test1
$ws.Cells.Item(64, 1).Value = "G1"
$ws.Cells.Item(64, 2).Value = "Test1"
$ws.Cells.Item(64, 3).Value = 45892
$ws.Cells.Item(64, 3).NumberFormat = $dateFormat
$ws.Cells.Item(64, 4).Value = 0.741922917787124
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = -0.01

# Row 65: G2 / sedrftgyhuioygtfrd
$ws.Cells.Item(65, 1).Value = "G2"
$ws.Cells.Item(65, 2).Value = "sedrftgyhuioygtfrd"
$ws.Cells.Item(65, 3).Value = 45892
$ws.Cells.Item(65, 3).NumberFormat = $dateFormat
$ws.Cells.Item(65, 4).Value = 0.741922917787124
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = -0.01

Write-Host "Added rows 64-65"
